# "Lopetasin kodutoo + excel"
# Update the "Nadal 3" (week 3) time log sheet: adjust a few end times and
# minute durations for rows 10, 12 and 13, and move the active-cell
# selection from F14 to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: end time moved from 22:30 to 23:00, duration 60 -> 90 minutes
$ws.Range("D10").Value = 0.95833333333333337
$ws.Range("F10").Value = 90

# Row 12: start time moved from 13:50 to 11:50, duration 60 -> 180 minutes
$ws.Range("C12").Value = 0.49305555555555558
$ws.Range("F12").Value = 180

# Row 13: end time moved from 19:30 to 18:30, duration 230 -> 170 minutes
$ws.Range("D13").Value = 0.77083333333333337
$ws.Range("F13").Value = 170

# F19 = SUM(F7:F18) recalculates automatically (700 -> 790)

# Move the selected/active cell from F14 to F11
$ws.Range("F11").Select()
